$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$rng = $ws.Range("D19:H19")
$rng.Hyperlinks.Delete()
$rng.ClearContents()
